# daily auto push: 2026-01-14 03:53 UTC
# Insert a new data row at row 647 (A1:D688 -> A1:D689), pushing the existing
# rows 647-688 down to 648-689, and populate the new row with:
#   A647 = 2026/01/14 (text)  B647 = 水  C647 = 11  D647 = 201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 647:688 down one row, leaving row 647 blank and ready for the
# new entry (mirrors Excel's "Insert Copied/Blank Cells -> Entire Row").
$ws.Rows.Item(647).Insert()

# Column B/C/D are plain text/number and round-trip through .Value cleanly.
$ws.Range("B647").Value = "水"
$ws.Range("C647").Value = 11
$ws.Range("D647").Value = 201

# Column A holds a date formatted as plain text ("2026/01/14"), like every
# other row in the sheet (cells are inline/shared strings, not real dates).
# A bare .Value assignment of a date-shaped string gets auto-converted to a
# date serial by Excel, so force literal text entry with a leading
# apostrophe (the standard "treat as text" trick) ...
$ws.Range("A647").Value = "'2026/01/14"

# ... then copy the (unstyled) number format from the sibling cell we just
# wrote so A647 ends up with the same default styling as the rest of the
# sheet instead of the "quote prefix" style that the apostrophe trick adds.
$ws.Range("B647").Copy()
$ws.Range("A647").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Inserted row 647: 2026/01/14, 水, 11, 201"
